# Tutorial 6 - updated attendance solution
# Updates the date formatting (slashes -> dashes) and recalculated
# attendance tallies (D=Total, E=Real, F=Duplicate, G=Invalid, H=Absent)
# for rows 3-21 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 11; Date = "25-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 12; Date = "29-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 13; Date = "01-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 14; Date = "05-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 15; Date = "08-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 16; Date = "12-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 2; E = 1; F = 1; G = 0; H = 0 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($item in $rows) {
    $r = $item.Row
    # Force the date column to remain plain text (matches the source
    # workbook's inline-string dates) instead of being auto-parsed
    # into a date serial number by Excel's smart input.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}
